$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Unhide rows 24-111 (previously hidden by the autofilter)
$ws.Rows("24:111").Hidden = $false

# 2. Clear the AutoFilter criteria on column D but keep the filter range/button
$ws.Range("D1:D114").AutoFilter(1)

# 3. Hide columns C, E, F, G, H:K, P, Q, S
$ws.Columns("C").Hidden = $true
$ws.Columns("E").Hidden = $true
$ws.Columns("F").Hidden = $true
$ws.Columns("G").Hidden = $true
$ws.Columns("H:K").Hidden = $true
$ws.Columns("P").Hidden = $true
$ws.Columns("Q").Hidden = $true
$ws.Columns("S").Hidden = $true

# 4. Update the selection to A1:R45 with active cell R45
$ws.Range("A1:R45").Select()
$ws.Range("R45").Activate()
